$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new blank rows at the correct alphabetically-sorted positions
#    (ascending row order so earlier inserts correctly push later rows down).
# ---------------------------------------------------------------------------
$ws.Rows("7:7").Insert()
$ws.Rows("9:9").Insert()
$ws.Rows("15:15").Insert()
$ws.Rows("17:17").Insert()

# ---------------------------------------------------------------------------
# 2. Re-create the merged-cell layout for each newly inserted row
#    (Insert() does not automatically re-merge the new blank row).
# ---------------------------------------------------------------------------
foreach ($r in @(7, 9, 15, 17)) {
    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

# ---------------------------------------------------------------------------
# 3. Write out the full product table (rows 7-22) exactly as it must read
#    after the edit - row number, name, balance, threshold, price, sell
#    price, transaction count.
# ---------------------------------------------------------------------------
$products = @(
    @(1,  "ABIMOL 300MG 5 RECTAL SUPP.",               "6:0",  "1", "15.00",  "15.0000", "1:0"),
    @(2,  "CEFOTAX 1GM I.M./I.V. VIAL - EIPICO",        "10:0", "1", "43.00",  "43.0000", "1:0"),
    @(3,  "CETAL 250MG/5ML 60ML SUSP",                  "19:0", "1", "31.00",  "31.0000", "1:0"),
    @(4,  "CONCOR COR 2.5MG 30 TABLETS",                "1:2",  "1", "60.00",  "60.0000", "1:0"),
    @(5,  "CONTROLOC 40MG 14 GASTRORESISTANT TAB",      "0:0",  "1", "188.00", "188.0000","1:0"),
    @(6,  "FORTAZEDIM 1 GM VIAL",                       "2:0",  "1", "59.00",  "118.0000","2:0"),
    @(7,  "GAST-REG 200 MG 30 TABS.",                   "1:2",  "1", "84.00",  "27.7200", "0:1"),
    @(8,  "MEBO 0.25% HERBAL AND NATURAL OINT. 15 GM",  "0:0",  "1", "109.00", "109.0000","1:0"),
    @(9,  "OPLEX-N SYRUP 125ML",                        "4:0",  "1", "31.00",  "31.0000", "1:0"),
    @(10, "SPASMOFEN 3 AMP. FOR I.M. INJ.",             "0:2",  "1", "39.00",  "12.8700", "0:1"),
    @(11, "TEGRETOL CR 400MG 20 F.C. DIVITABS",         "0:1",  "1", "106.00", "53.0000", "0:1"),
    @(12, "TICANASE 0.05% NASAL SPRAY 12 GM",           "1:0",  "1", "70.00",  "70.0000", "1:0"),
    @(13, "ZURCAL 40MG 14 GASTRO RESISTANT TAB",        "6:0",  "1", "96.00",  "96.0000", "1:0"),
    @(14, "بلاستر مترسيلك 2 سم",                        "19:0", "0", "15.00",  "15.0000", "1:0"),
    @(15, "سرنجات 3 سم",                                 "0:0",  "0", "2.00",   "12.0000", "6:0"),
    @(16, "كالونا ",                                     "0:0",  "0", "15.00",  "15.0000", "1:0")
)

$row = 7
foreach ($p in $products) {
    $ws.Range("A$row").Value = $p[0]
    $ws.Range("C$row").Value = $p[1]
    $ws.Range("H$row").Value = $p[2]
    $ws.Range("L$row").Value = $p[3]
    $ws.Range("N$row").Value = $p[4]
    $ws.Range("P$row").Value = $p[5]
    $ws.Range("Q$row").Value = $p[6]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4. Update the grand-total cell (now on row 23) and the footer timestamp
#    (now on row 24).
# ---------------------------------------------------------------------------
$ws.Range("P23").Value = 896.59000000000003
$ws.Range("A24").Value = "Thursday, 29 May, 2025 12:28 PM"
